$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Header Bar Hidding Content" bug (row 3) as Fixed, matching
# the existing "Yes" value used for the other fixed bug in row 2 (F2).
$ws.Range("F3").Value = "Yes"
